$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Version 6.0.0 -> 6.1.0
$wsMetadata.Range("B3").Value = "6.1.0"

# Date 2022-01-21T20:46:54+00:00 -> 2022-05-31T20:10:14+00:00
$wsMetadata.Range("B8").Value = "2022-05-31T20:10:14+00:00"

# Definitions for Engagement Attributed PCP / Specialist
$wsConcepts.Range("D4").Value = "Primary care physician attributed by the patient engagement pipeline"
$wsConcepts.Range("D5").Value = "Specialist physician attributed by the patient engagement pipeline"

$wb.Save()
